$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-pulled the data: two new columns ("height" and "weight")
# are inserted right before the existing "fantasy points" column, which
# slides from E to G. New layout:
#   B=rec_yds  C=rec_td  D=fumbles  E=height  F=weight  G=fantasy points

$lastRow = 16

# 1) Shift the existing "fantasy points" values from column E to the new
#    column G (bottom-up so we read each value before anything overwrites
#    it).
for ($r = $lastRow; $r -ge 2; $r--) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
}

# 2) Fill in the new height/weight values for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 260
}

# 3) Headers. E1 already holds a shared string ("fantasy points"); relabel
#    it "height" and add the two new headers, then match the header
#    formatting (bold, centered, bordered) already used on B1:E1.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
